# Generate Report for Handoff
#
# The localization-status report records, per localized file, the most
# recent timestamp at which a handoff (xliff) package was generated.
# The "37d512c1-5525-4828-a0b9-217ee8520cf5.md" file has just had a new
# handoff package generated for the de-de locale, so both:
#   - the per-locale worksheet ("de-de") "Latest Handoff Datetime" cell, and
#   - the "Overview" worksheet's roll-up "Latest HO Xliff Generate Date" cell
# need to be updated to the new generation timestamp.

$wb = $excel.ActiveWorkbook

$newHandoffDatetime = "2016-08-31 06:46:48"

$overview = $wb.Worksheets.Item("Overview")
$dede = $wb.Worksheets.Item("de-de")

# Locate the row for the file that was just handed off, rather than
# assuming a fixed row number, so the script is resilient to row order.
$targetFile = "37d512c1-5525-4828-a0b9-217ee8520cf5.md"

$overviewRows = $overview.UsedRange.Rows.Count
for ($r = 2; $r -le $overviewRows; $r++) {
    $cellValue = $overview.Cells.Item($r, 1).Value2
    if ($cellValue -eq $targetFile) {
        $overview.Cells.Item($r, 7).Value2 = $newHandoffDatetime
    }
}

$dedeRows = $dede.UsedRange.Rows.Count
for ($r = 2; $r -le $dedeRows; $r++) {
    $cellValue = $dede.Cells.Item($r, 1).Value2
    if ($cellValue -eq $targetFile) {
        $dede.Cells.Item($r, 8).Value2 = $newHandoffDatetime
    }
}
